$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44503

# Row 3
$ws.Range("D3").Value = 44510
$ws.Range("J3").Value = 600
$ws.Range("K3").Value = 900
$ws.Range("L3").Value = 1000
$ws.Range("M3").Value = 950
$ws.Range("P3").Value = 950

# Row 4
$ws.Range("D4").Value = 44511
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 900
$ws.Range("L4").Value = 1000
$ws.Range("M4").Value = 950
$ws.Range("P4").Value = 950

# Row 5
$ws.Range("D5").Value = 44512
$ws.Range("J5").Value = 600
$ws.Range("K5").Value = 900
$ws.Range("L5").Value = 1000
$ws.Range("M5").Value = 950
$ws.Range("P5").Value = 950

# Row 6
$ws.Range("D6").Value = 44504
$ws.Range("J6").Value = 500

# Row 7
$ws.Range("D7").Value = 44508
$ws.Range("J7").Value = 400

# Row 8
$ws.Range("D8").Value = 44525
$ws.Range("J8").Value = 360

# Row 9
$ws.Range("D9").Value = 44523
$ws.Range("J9").Value = 400
$ws.Range("K9").Value = 800
$ws.Range("L9").Value = 900
$ws.Range("M9").Value = 850
$ws.Range("P9").Value = 850

# Row 10
$ws.Range("D10").Value = 44530
$ws.Range("J10").Value = 300
$ws.Range("K10").Value = 800
$ws.Range("L10").Value = 900
$ws.Range("M10").Value = 850
$ws.Range("P10").Value = 850

# Row 11
$ws.Range("D11").Value = 44517
$ws.Range("J11").Value = 500

# Row 12
$ws.Range("D12").Value = 44532
$ws.Range("J12").Value = 240

# Row 13
$ws.Range("D13").Value = 44537
$ws.Range("J13").Value = 400

# Row 16
$ws.Range("D16").Value = 44545
$ws.Range("J16").Value = 4000

# Row 17
$ws.Range("D17").Value = 44553
$ws.Range("J17").Value = 8000

# Row 18
$ws.Range("D18").Value = 44524

# Row 19
$ws.Range("D19").Value = 44518
$ws.Range("J19").Value = 400
$ws.Range("K19").Value = 800
$ws.Range("L19").Value = 900
$ws.Range("M19").Value = 850
$ws.Range("P19").Value = 850

# Row 20
$ws.Range("D20").Value = 44505
$ws.Range("J20").Value = 440
